$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 8, shifting existing rows 8:104 down to 9:105
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new data record
$ws.Cells.Item(8, 1).Value = 8
$ws.Cells.Item(8, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = Get-Date -Year 2021 -Month 12 -Day 15 -Hour 0 -Minute 0 -Second 0 -Millisecond 0
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100112001
$ws.Cells.Item(8, 7).Value = "Berenjena"
$ws.Cells.Item(8, 8).Value = "Sin especificar"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 540
$ws.Cells.Item(8, 11).Value = 10000
$ws.Cells.Item(8, 12).Value = 11000
$ws.Cells.Item(8, 13).Value = 10500
$ws.Cells.Item(8, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(8, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(8, 16).Value = 175
$ws.Cells.Item(8, 17).Value = 60
$ws.Cells.Item(8, 18).Value = "Hortaliza"
